# Populate the remaining header cells in row 1 (D1:L1) with their column
# letters, matching the expanded test-data header row, and restore the
# selection to C1 (top of the data-fetch range) instead of the stale G33.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "D"
$ws.Range("E1").Value = "E"
$ws.Range("F1").Value = "F"
$ws.Range("G1").Value = "G"
$ws.Range("H1").Value = "H"
$ws.Range("I1").Value = "I"
$ws.Range("J1").Value = "J"
$ws.Range("K1").Value = "K"
$ws.Range("L1").Value = "L"

$ws.Range("C1").Select()
